$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: fill in end time, and the elapsed-time formula (matching the
# existing D-column shared formula pattern/format used by D7:D16).
$ws.Range("C17").Value = 0.875
$ws.Range("D17").Formula = "=C17-B17"
$ws.Range("D17").NumberFormat = "h:mm:ss;@"

# Row 18: new time-log entry (Test Plan Document work).
$ws.Range("B18").Value = 0.875
$ws.Range("C18").Value = 0.91666666666666663
$ws.Range("D18").Formula = "=C18-B18"
$ws.Range("D18").NumberFormat = "h:mm:ss;@"
$ws.Range("E18").Value = "Test Plan Document"
$ws.Range("F18").Value = "Creating Test Plan Document"

# New formatted-but-empty cell further down the sheet (extends the used
# range / dimension to F33) with an h:mm number format.
$ws.Range("D33").NumberFormat = "h:mm"

# Restore the active selection to where the user last left off editing.
$ws.Range("E18").Select()
